$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 5.207399999999996
$ws.Range("E4").Value = 13.8224
$ws.Range("E5").Value = 13.13369999999999
$ws.Range("B6").Value = 9.276499999999992
$ws.Range("E6").Value = 12.5066
$ws.Range("B7").Value = 5.0406
$ws.Range("B8").Value = 5.050999999999993
$ws.Range("E8").Value = 14.2064
$ws.Range("B16").Value = 9.131300000000007
$ws.Range("E16").Value = 12.3603
$ws.Range("B20").Value = 5.506999999999997
$ws.Range("B21").Value = 5.122800000000001
$ws.Range("E22").Value = 12.13389999999999
